# ECE 383 Lec 9 - update the title slide's subtitle line from
# "Lecture 7 - Finite State Machines" to "Lecture 9 - Finite State Machines",
# splitting the run so the new lesson number ("9") lives in its own run
# (mirrors the way PowerPoint splits a run when the user edits only part
# of it), instead of just doing a whole-string replace.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)          # "Rectangle 2" - the ctrTitle placeholder
$tr = $sh.TextFrame.TextRange

$full = $tr.Text

# Locate "Lecture 7" inside the title text rather than hard-coding character
# offsets, so the edit still lands correctly even if earlier text shifts.
$lecIdx = $full.IndexOf("Lecture 7")
if ($lecIdx -lt 0) {
    throw "Could not find 'Lecture 7' in the title shape text"
}

# TextRange.Characters uses 1-based character positions.
$lecStart = $lecIdx + 1            # start of "Lecture 7 - Finite..."
$digitPos = $lecStart + 8          # "Lecture " is 8 characters -> the "7"
$spacePos = $digitPos + 1          # the space right after the digit

# 1) "7" -> "9" : this splits the original run into
#       "Lecture " | "9" | " - Finite State Machines"
$digitRange = $tr.Characters($digitPos, 1)
$digitRange.Text = "9"

# 2) Split the leading space off of " - Finite State Machines" so the dash
#    clause becomes its own run:
#       "Lecture " | "9" | " " | "- Finite State Machines"
$spaceRange = $tr.Characters($spacePos, 1)
$spaceRange.Text = " "

Write-Output ("Title is now: " + $tr.Text)
